$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 7
$ws.Range("D7").Value = 44431
$ws.Range("N7").Value = 30000
$ws.Range("O7").Value = 30000
$ws.Range("P7").Value = 30000
$ws.Range("S7").Value = 3000
# Row 8
$ws.Range("D8").Value = 44459
$ws.Range("M8").Value = 100
# Row 9
$ws.Range("D9").Value = 44435
$ws.Range("M9").Value = 160
# Row 10
$ws.Range("D10").Value = 44445
$ws.Range("L10").Value = 'Primera'
$ws.Range("M10").Value = 250
$ws.Range("N10").Value = 28000
$ws.Range("P10").Value = 29200
$ws.Range("S10").Value = 2920
# Row 11
$ws.Range("D11").Value = 44446
$ws.Range("M11").Value = 200
$ws.Range("N11").Value = 30000
$ws.Range("P11").Value = 30000
$ws.Range("S11").Value = 3000
# Row 12
$ws.Range("D12").Value = 44467
$ws.Range("L12").Value = 'Especial'
$ws.Range("M12").Value = 100
# Row 13
$ws.Range("D13").Value = 44467
$ws.Range("L13").Value = 'Primera'
$ws.Range("M13").Value = 100
$ws.Range("N13").Value = 28000
$ws.Range("O13").Value = 28000
$ws.Range("P13").Value = 28000
$ws.Range("S13").Value = 2800
# Row 14
$ws.Range("D14").Value = 44434
$ws.Range("L14").Value = 'Especial'
$ws.Range("M14").Value = 60
# Row 15
$ws.Range("D15").Value = 44441
$ws.Range("L15").Value = 'Primera'
$ws.Range("M15").Value = 150
# Row 16
$ws.Range("D16").Value = 44453
$ws.Range("M16").Value = 135
# Row 17
$ws.Range("L17").Value = 'Especial'
$ws.Range("M17").Value = 100
$ws.Range("N17").Value = 30000
$ws.Range("O17").Value = 30000
$ws.Range("P17").Value = 30000
$ws.Range("S17").Value = 3000
# Row 18
$ws.Range("D18").Value = 44448
$ws.Range("L18").Value = 'Primera'
$ws.Range("M18").Value = 80
$ws.Range("N18").Value = 28000
$ws.Range("O18").Value = 28000
$ws.Range("P18").Value = 28000
$ws.Range("S18").Value = 2800
# Row 19
$ws.Range("D19").Value = 44455
$ws.Range("M19").Value = 150
# Row 20
$ws.Range("D20").Value = 44460
$ws.Range("M20").Value = 80
# Row 21
$ws.Range("D21").Value = 44466
$ws.Range("M21").Value = 110
$ws.Range("N21").Value = 30000
$ws.Range("O21").Value = 30000
$ws.Range("P21").Value = 30000
$ws.Range("S21").Value = 3000
# Row 22
$ws.Range("D22").Value = 44447
$ws.Range("M22").Value = 50
$ws.Range("N22").Value = 32000
$ws.Range("O22").Value = 32000
$ws.Range("P22").Value = 32000
$ws.Range("S22").Value = 3200
# Row 23
$ws.Range("L23").Value = 'Especial'
$ws.Range("M23").Value = 150
$ws.Range("N23").Value = 30000
$ws.Range("O23").Value = 30000
$ws.Range("P23").Value = 30000
$ws.Range("S23").Value = 3000
# Row 24
$ws.Range("D24").Value = 44463
$ws.Range("N24").Value = 26000
$ws.Range("O24").Value = 26000
$ws.Range("P24").Value = 26000
$ws.Range("S24").Value = 2600
# Row 25
$ws.Range("D25").Value = 44438
$ws.Range("L25").Value = 'Primera'
$ws.Range("M25").Value = 100
# Row 26
$ws.Range("L26").Value = 'Especial'
$ws.Range("M26").Value = 150
$ws.Range("N26").Value = 30000
$ws.Range("O26").Value = 30000
$ws.Range("P26").Value = 30000
$ws.Range("S26").Value = 3000
# Row 27
$ws.Range("A27").Value = 5
$ws.Range("B27").Value = 'Macroferia Regional de Talca'
$ws.Range("C27").Value = 'Maule'
$ws.Range("D27").Value = 44461
$ws.Range("E27").Value = 7
$ws.Range("F27").Value = 'Fruta'
$ws.Range("G27").Value = 100107
$ws.Range("H27").Value = 'Otros'
$ws.Range("I27").Value = 100107002
$ws.Range("J27").Value = 'Chirimoya'
$ws.Range("K27").Value = 'Cultivar IV Región'
$ws.Range("L27").Value = 'Primera'
$ws.Range("M27").Value = 100
$ws.Range("N27").Value = 25000
$ws.Range("O27").Value = 25000
$ws.Range("P27").Value = 25000
$ws.Range("Q27").Value = '$/bandeja 10 kilos'
$ws.Range("R27").Value = 'Provincia de Limarí'
$ws.Range("S27").Value = 2500
$ws.Range("T27").Value = 10

# Ensure date formatting for the new row 27 date cell (new cells default to General)
$ws.Range("D27").NumberFormat = "YYYY-MM-DD HH:MM:SS"

Write-Host "Edit applied"
